$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Delete column D (oldest quarter, 1399/06) - this shifts everything left by one column
$ws.Columns.Item(4).Delete()

# New newest-quarter values for column M (previously column N before the shift; after the delete of D,
# the old column M becomes L, and we fill the new, 11th quarter data into column M)
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# "1402-02-25" looks like a date to Excel's parser, so write it as a text
# formula first and then flatten it to a static value to keep it a plain
# string (and keep the existing cell style/format untouched).
$ws.Range("M9").Formula = "=""1402-02-25"""
$ws.Range("M9").Copy()
$ws.Range("M9").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("M11").Value = 1702125
$ws.Range("M12").Value = -1530236
$ws.Range("M13").Value = 171888
$ws.Range("M14").Value = -12910
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = 752
$ws.Range("M17").Value = 159730
$ws.Range("M18").Value = -305
$ws.Range("M19").Value = 12346
$ws.Range("M20").Value = 171771
$ws.Range("M21").Value = 42540
$ws.Range("M22").Value = 214311
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 214311
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 295562
$ws.Range("M27").Value = 0

# Copy formatting from column L into column M for consistency
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)  # xlPasteFormats
